# ------------------------------------------------------------------
# Adds a new "2022-Q3" sheet of fund-holding data and records a new
# summary row for it on the "总计" (totals) sheet.
# ------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# --------------------------------------------------------------
# 1) Update the "总计" summary sheet: insert a new row 2 for 2022-Q3
#    and shift the existing quarters (and their running index in
#    column A) down by one row.
# --------------------------------------------------------------
$summary = $wb.Worksheets.Item(1)

$summary.Rows.Item(2).Insert()

# Copy formatting (bold/border style on col A) down from the row
# that used to be row 2 (now row 3) onto the freshly inserted row.
$summary.Range("A3:D3").Copy()
$summary.Range("A2:D2").PasteSpecial(-4122)   # xlPasteFormats

$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q3"
$summary.Range("C2").Value = 6
$summary.Range("D2").Value = 0.16

$summary.Range("A3").Value = 1
$summary.Range("A4").Value = 2
$summary.Range("A5").Value = 3
$summary.Range("A6").Value = 4
$summary.Range("A7").Value = 5

# --------------------------------------------------------------
# 2) Insert a brand-new "2022-Q3" worksheet, positioned before the
#    existing "2022-Q2" tab, populated with the quarter's fund
#    holding data. Duplicating the "2022-Q2" sheet first gives us
#    an identical layout/format (header row, column widths, etc.)
#    to build on.
# --------------------------------------------------------------
$q2sheet = $wb.Worksheets.Item("2022-Q2")
$q2sheet.Copy($q2sheet)
$newSheet = $wb.Worksheets.Item("2022-Q2 (2)")
$newSheet.Name = "2022-Q3"

# The template (2022-Q2) only has 4 data rows (rows 2-5); 2022-Q3
# needs 6 (rows 2-7), so insert two more rows, copying the format
# of the last template row down onto them.
$newSheet.Rows.Item(6).Insert()
$newSheet.Range("A5:H5").Copy()
$newSheet.Range("A6:H6").PasteSpecial(-4122)   # xlPasteFormats

$newSheet.Rows.Item(7).Insert()
$newSheet.Range("A5:H5").Copy()
$newSheet.Range("A7:H7").PasteSpecial(-4122)   # xlPasteFormats

# Columns B:G hold text-like values (fund codes, names and numbers
# that must keep their exact textual representation, e.g. leading
# zeros, trailing zeros after a decimal point, etc). Prefixing the
# value with a leading apostrophe forces Excel to store it as Text
# instead of silently converting numeric-looking strings to real
# numbers (which would drop formatting like "080005" -> 80005).
$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").Value = "'080005"
$newSheet.Range("C2").Value = "长盛量化红利混合"
$newSheet.Range("D2").Value = "'1.89"
$newSheet.Range("E2").Value = "'61.68"
$newSheet.Range("F2").Value = "'2.12"
$newSheet.Range("G2").Value = "'0.0401"
$newSheet.Range("H2").Value = 9

$newSheet.Range("A3").Value = 1
$newSheet.Range("B3").Value = "'000573"
$newSheet.Range("C3").Value = "天弘通利混合"
$newSheet.Range("D3").Value = "'1.01"
$newSheet.Range("E3").Value = "'79.25"
$newSheet.Range("F3").Value = "'3.02"
$newSheet.Range("G3").Value = "'0.0305"
$newSheet.Range("H3").Value = 8

$newSheet.Range("A4").Value = 2
$newSheet.Range("B4").Value = "'006700"
$newSheet.Range("C4").Value = "红土创新稳健混合A"
$newSheet.Range("D4").Value = "'0.66"
$newSheet.Range("E4").Value = "'27.52"
$newSheet.Range("F4").Value = "'4.53"
$newSheet.Range("G4").Value = "'0.0299"
$newSheet.Range("H4").Value = 3

$newSheet.Range("A5").Value = 3
$newSheet.Range("B5").Value = "'006701"
$newSheet.Range("C5").Value = "红土创新稳健混合C"
$newSheet.Range("D5").Value = "'0.50"
$newSheet.Range("E5").Value = "'27.52"
$newSheet.Range("F5").Value = "'4.53"
$newSheet.Range("G5").Value = "'0.0226"
$newSheet.Range("H5").Value = 3

$newSheet.Range("A6").Value = 4
$newSheet.Range("B6").Value = "'011198"
$newSheet.Range("C6").Value = "交银施罗德鑫选回报混合A"
$newSheet.Range("D6").Value = "'4.49"
$newSheet.Range("E6").Value = "'20.45"
$newSheet.Range("F6").Value = "'0.47"
$newSheet.Range("G6").Value = "'0.0211"
$newSheet.Range("H6").Value = 8

$newSheet.Range("A7").Value = 5
$newSheet.Range("B7").Value = "'011199"
$newSheet.Range("C7").Value = "交银施罗德鑫选回报混合C"
$newSheet.Range("D7").Value = "'3.20"
$newSheet.Range("E7").Value = "'20.45"
$newSheet.Range("F7").Value = "'0.47"
$newSheet.Range("G7").Value = "'0.0150"
$newSheet.Range("H7").Value = 8

# Clear the "quote prefix" formatting flag left behind by the
# leading apostrophes so the cells end up with the default
# (unstyled) look, matching the rest of the sheet.
$newSheet.Range("B2:G7").ClearFormats()
